$p = $ppt.ActivePresentation

# --- 1. Update every "2020/4/1" date placeholder field to "2020/6/27" --------
# The decks re-cache the datetime{...} fields on save; every slide that still
# shows the stale "2020/4/1" text gets the new date.
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $sh = $s.Shapes.Item($i)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq "2020/4/1") {
                $sh.TextFrame.TextRange.Text = "2020/6/27"
            }
        }
    }
}

# --- 2. "More on Why React?" slide: shrink the bullet box and drop the ------
#        trailing empty bullet paragraph at the end of the list.
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $sh = $s.Shapes.Item($i)
        if ($sh.HasTextFrame) {
            $txt = $sh.TextFrame.TextRange.Text
            if ($txt.StartsWith("More on Why React?")) {
                # Resize: keep position/width, only shorten the height.
                $sh.Height = 3600401 / 12700

                # Remove the last (empty) paragraph of the bulleted list.
                $tr = $sh.TextFrame.TextRange
                $paraCount = $tr.Paragraphs().Count
                $lastPara = $tr.Paragraphs($paraCount, 1)
                if ($lastPara.Text -eq "") {
                    $lastPara.Delete()
                }
            }
        }
    }
}
